$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1575
$ws.Range("I8").Value = 150
$ws.Range("K8").Value = 450
$ws.Range("M8").Value = -311
$ws.Range("H40").Value = 2486.75
$ws.Range("I40").Value = 1449.5
$ws.Range("J40").Value = 2832.5
$ws.Range("K40").Value = 1449.5
$ws.Range("L40").Value = 2832.5
$ws.Range("M40").Value = -1274.5
$ws.Range("N40").Value = -3182.5
$ws.Range("H86").Value = 14361328
$ws.Range("I86").Value = 4831.6665
$ws.Range("K86").Value = 4831.6665
$ws.Range("M86").Value = -3708.6665
$ws.Range("H88").Value = 1492
$ws.Range("I88").Value = 1326.3334
$ws.Range("J88").Value = 1740.5
$ws.Range("K88").Value = 1326.3334
$ws.Range("L88").Value = 1740.5
$ws.Range("M88").Value = -920.3334
$ws.Range("N88").Value = -2552.5
$ws.Range("H89").Value = 14361328
$ws.Range("I89").Value = 4831.6665
$ws.Range("K89").Value = 24158.3325
$ws.Range("M89").Value = -18542.3325
$ws.Range("H91").Value = 1492
$ws.Range("I91").Value = 1326.3334
$ws.Range("J91").Value = 1740.5
$ws.Range("K91").Value = 1326.3334
$ws.Range("L91").Value = 1740.5
$ws.Range("M91").Value = 77.66660000000002
$ws.Range("N91").Value = -4548.5
$ws.Range("H92").Value = 66926.734
$ws.Range("I92").Value = 278.64285
$ws.Range("K92").Value = 278.64285
$ws.Range("M92").Value = 969.35715
$ws.Range("H111").Value = 29647.889
$ws.Range("I111").Value = 6681.143
$ws.Range("J111").Value = 110031.5
$ws.Range("K111").Value = 20043.429
$ws.Range("L111").Value = 330094.5
$ws.Range("M111").Value = -16976.429
$ws.Range("N111").Value = -336228.5
$ws.Range("H137").Value = 40001230
$ws.Range("I137").Value = 55556840
$ws.Range("J137").Value = 1099.2858
$ws.Range("K137").Value = 166670520
$ws.Range("L137").Value = 3297.8574
$ws.Range("M137").Value = -166667970
$ws.Range("N137").Value = -8397.857400000001
$ws.Range("H138").Value = 1732.0526
$ws.Range("I138").Value = 1348.3572
$ws.Range("K138").Value = 4045.0716
$ws.Range("M138").Value = 1094.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 3449.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 3449.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 3449.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -3737.5
$ws.Range("H61").Value = 1303.8182
$ws.Range("I61").Value = 1136.6
$ws.Range("K61").Value = 1136.6
$ws.Range("M61").Value = -924.5999999999999
$ws.Range("H88").Value = 18521450
$ws.Range("I88").Value = 27779760
$ws.Range("K88").Value = 27779760
$ws.Range("M88").Value = -27779354
$ws.Range("H91").Value = 18521450
$ws.Range("I91").Value = 27779760
$ws.Range("K91").Value = 27779760
$ws.Range("M91").Value = -27778356
$ws.Range("H122").Value = 1359.8788
$ws.Range("I122").Value = 1082.24
$ws.Range("K122").Value = 3246.72
$ws.Range("M122").Value = -796.7200000000003
$ws.Range("H132").Value = 915.0540999999999
$ws.Range("I132").Value = 846.84375
$ws.Range("K132").Value = 2540.53125
$ws.Range("M132").Value = -10.53125
$ws.Range("H135").Value = 73607.25
$ws.Range("J135").Value = 73607.25
$ws.Range("L135").Value = 73607.25
$ws.Range("N135").Value = -83747.25
$ws.Range("H136").Value = 1303.8182
$ws.Range("I136").Value = 1136.6
$ws.Range("K136").Value = 3409.8
$ws.Range("M136").Value = -859.7999999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 49800
$ws.Range("J135").Value = 49800
$ws.Range("L135").Value = 49800
$ws.Range("N135").Value = -59940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 937.8333
$ws.Range("J5").Value = 1275.5555
$ws.Range("L5").Value = 3826.6665
$ws.Range("N5").Value = -4050.6665
$ws.Range("H6").Value = 655.5
$ws.Range("I6").Value = 78
$ws.Range("J6").Value = 1002
$ws.Range("K6").Value = 234
$ws.Range("L6").Value = 3006
$ws.Range("M6").Value = -121
$ws.Range("N6").Value = -3232
$ws.Range("H56").Value = 6043.8
$ws.Range("I56").Value = 6043.8
$ws.Range("K56").Value = 6043.8
$ws.Range("M56").Value = -5513.8
$ws.Range("H132").Value = 2069
$ws.Range("J132").Value = 2210
$ws.Range("L132").Value = 19890
$ws.Range("N132").Value = -24950
$ws.Range("H135").Value = 937.8333
$ws.Range("J135").Value = 1275.5555
$ws.Range("L135").Value = 11479.9995
$ws.Range("N135").Value = -16549.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1340.6666
$ws.Range("I9").Value = 431
$ws.Range("J9").Value = 1522.6
$ws.Range("K9").Value = 431
$ws.Range("L9").Value = 1522.6
$ws.Range("M9").Value = -261
$ws.Range("N9").Value = -1862.6
$ws.Range("H11").Value = 7186177.5
$ws.Range("I11").Value = 5156473.5
$ws.Range("J11").Value = 8404000
$ws.Range("K11").Value = 5156473.5
$ws.Range("L11").Value = 8404000
$ws.Range("M11").Value = -5156334.5
$ws.Range("N11").Value = -8404278
$ws.Range("H69").Value = 45000
$ws.Range("J69").Value = 45000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46498
$ws.Range("H72").Value = 45000
$ws.Range("J72").Value = 45000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -142488
$ws.Range("H80").Value = 3437.3572
$ws.Range("I80").Value = 3269.9443
$ws.Range("J80").Value = 3738.7
$ws.Range("K80").Value = 3269.9443
$ws.Range("L80").Value = 3738.7
$ws.Range("M80").Value = -2271.9443
$ws.Range("N80").Value = -5734.7
$ws.Range("H83").Value = 3437.3572
$ws.Range("I83").Value = 3269.9443
$ws.Range("J83").Value = 3738.7
$ws.Range("K83").Value = 16349.7215
$ws.Range("L83").Value = 18693.5
$ws.Range("M83").Value = -11357.7215
$ws.Range("N83").Value = -28677.5
$ws.Range("H102").Value = 1387.9286
$ws.Range("I102").Value = 1065.591
$ws.Range("K102").Value = 1065.591
$ws.Range("M102").Value = 556.4090000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 906.63635
$ws.Range("I16").Value = 597.4
$ws.Range("J16").Value = 3999
$ws.Range("K16").Value = 597.4
$ws.Range("L16").Value = 3999
$ws.Range("M16").Value = -427.4
$ws.Range("N16").Value = -4339
$ws.Range("H19").Value = 897
$ws.Range("I19").Value = 897
$ws.Range("K19").Value = 897
$ws.Range("M19").Value = -727
$ws.Range("H22").Value = 309.16666
$ws.Range("I22").Value = 171.2
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 171.2
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = 123.8
$ws.Range("N22").Value = -1589
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H27").Value = 309.16666
$ws.Range("I27").Value = 171.2
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 171.2
$ws.Range("L27").Value = 999
$ws.Range("M27").Value = -64.19999999999999
$ws.Range("N27").Value = -1213
$ws.Range("H46").Value = 1913.7354
$ws.Range("I46").Value = 1329.4445
$ws.Range("J46").Value = 2571.0625
$ws.Range("K46").Value = 1329.4445
$ws.Range("L46").Value = 2571.0625
$ws.Range("M46").Value = -1141.4445
$ws.Range("N46").Value = -2947.0625
$ws.Range("H55").Value = 481.37036
$ws.Range("I55").Value = 336.31818
$ws.Range("K55").Value = 336.31818
$ws.Range("M55").Value = -163.31818
$ws.Range("H61").Value = 1144.1666
$ws.Range("I61").Value = 948.2727
$ws.Range("K61").Value = 948.2727
$ws.Range("M61").Value = -746.2727
$ws.Range("H113").Value = 1144.1666
$ws.Range("I113").Value = 948.2727
$ws.Range("K113").Value = 948.2727
$ws.Range("M113").Value = 1221.7273
$ws.Range("H122").Value = 3943.1482
$ws.Range("I122").Value = 2787.6316
$ws.Range("K122").Value = 8362.8948
$ws.Range("M122").Value = -5912.8948
$ws.Range("H133").Value = 272687.5
$ws.Range("J133").Value = 272687.5
$ws.Range("L133").Value = 272687.5
$ws.Range("N133").Value = -277747.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1408341.5
$ws.Range("J62").Value = 8673.5
$ws.Range("L62").Value = 8673.5
$ws.Range("N62").Value = -9921.5
$ws.Range("H65").Value = 1408341.5
$ws.Range("J65").Value = 8673.5
$ws.Range("L65").Value = 43367.5
$ws.Range("N65").Value = -49607.5
$ws.Range("H107").Value = 764.4286
$ws.Range("I107").Value = 749
$ws.Range("J107").Value = 770.6
$ws.Range("K107").Value = 2247
$ws.Range("L107").Value = 2311.8
$ws.Range("M107").Value = -327
$ws.Range("N107").Value = -6151.8
$ws.Range("H132").Value = 1497.3625
$ws.Range("I132").Value = 1326.9552
$ws.Range("K132").Value = 3980.8656
$ws.Range("M132").Value = -1450.8656
